# The product-import template's "Descripción" column header is being
# renamed to "Impresoras" (see sharedStrings.xml si index 2, used by C1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Impresoras"
